$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column width adjustments (cols B, C, E, F) ---
$ws.Columns.Item(2).ColumnWidth = 22.944010416666668
$ws.Columns.Item(3).ColumnWidth = 21.053385416666668
$ws.Columns.Item(5).ColumnWidth = 22.498697916666668
$ws.Columns.Item(6).ColumnWidth = 13.385416666666666

# --- New cell values (fills in print-flag lookup table) ---
$ws.Range('G30').Value = 'pribcf'
$ws.Range('G31').Value = 'pric'
$ws.Range('G32').Value = 'prikd'
$ws.Range('G34').Value = 'prigfb'
$ws.Range('G40').Value = 'prip'
$ws.Range('G43').Value = 'privel'
$ws.Range('G44').Value = 'priwel'
$ws.Range('G47').Value = 'primaphead'
$ws.Range('G48').Value = 'primapv'
$ws.Range('E30').Value = 'print_bc_flows'
$ws.Range('F30').Value = 'prbcf'
$ws.Range('G30').Value = 'pribcf'
$ws.Range('H30').Value = 'timprbcf'
$ws.Range('J30').Value = 'ntprbcf'
$ws.Range('E31').Value = 'print_components'
$ws.Range('F31').Value = 'prc'
$ws.Range('G31').Value = 'pric'
$ws.Range('H31').Value = 'timprc'
$ws.Range('J31').Value = 'ntprc'
$ws.Range('E32').Value = 'print_conductances'
$ws.Range('F32').Value = 'prkd'
$ws.Range('G32').Value = 'prikd'
$ws.Range('H32').Value = 'timprkd'
$ws.Range('J32').Value = 'ntprkd'
$ws.Range('E34').Value = 'print_global_flow_balance'
$ws.Range('F34').Value = 'prgfb'
$ws.Range('G34').Value = 'prigfb'
$ws.Range('H34').Value = 'timprgfb'
$ws.Range('J34').Value = 'ntprgfb'
$ws.Range('E35').Value = 'print_force_chemistry'
$ws.Range('F35').Value = 'prf_chem_phrq'
$ws.Range('G35').Value = 'priforce_chem_phrq'
$ws.Range('H35').Value = 'timprfchem'
$ws.Range('I35').Value = 'prf_chem_phrqi'
$ws.Range('E40').Value = 'print_heads'
$ws.Range('F40').Value = 'prp'
$ws.Range('G40').Value = 'prip'
$ws.Range('H40').Value = 'timprp'
$ws.Range('J40').Value = 'ntprp'
$ws.Range('E41').Value = 'print_restart_hst'
$ws.Range('F41').Value = 'prcpd'
$ws.Range('G41').Value = 'pricpd'
$ws.Range('H41').Value = 'timprcpd'
$ws.Range('J41').Value = 'ntprcpd'
$ws.Range('E43').Value = 'print_velocities'
$ws.Range('F43').Value = 'prvel'
$ws.Range('G43').Value = 'privel'
$ws.Range('H43').Value = 'timprvel'
$ws.Range('J43').Value = 'ntprvel'
$ws.Range('E44').Value = 'print_wells'
$ws.Range('F44').Value = 'prwel'
$ws.Range('G44').Value = 'priwel'
$ws.Range('H44').Value = 'timprwel'
$ws.Range('J44').Value = 'ntprwel'
$ws.Range('E46').Value = 'print_xyz_components'
$ws.Range('F46').Value = 'prmapc'
$ws.Range('G46').Value = 'primapcomp'
$ws.Range('H46').Value = 'timprmapc'
$ws.Range('J46').Value = 'ntprmapcomp'
$ws.Range('E47').Value = 'print_xyz_heads'
$ws.Range('F47').Value = 'prmaph'
$ws.Range('G47').Value = 'primaphead'
$ws.Range('H47').Value = 'timprmaph'
$ws.Range('J47').Value = 'ntprmaphead'
$ws.Range('E48').Value = 'print_xyz_velocities'
$ws.Range('F48').Value = 'vecmap/primapv'
$ws.Range('G48').Value = 'primapv'
$ws.Range('H48').Value = 'timprmapv'
$ws.Range('J48').Value = 'ntprmapv'
$ws.Range('E49').Value = 'print_xyz_wells'
$ws.Range('F49').Value = 'prtem'
$ws.Range('G49').Value = 'pri_well_timser'
$ws.Range('H49').Value = 'timprtem'
$ws.Range('J49').Value = 'ntprtem'
$ws.Range('E52').Value = 'print_zone_flows_tsv'
$ws.Range('F52').Value = 'przf_tsv'
$ws.Range('G52').Value = 'pri_zf_tsv'
$ws.Range('H52').Value = 'timprzf_tsv'
$ws.Range('J52').Value = 'ntprzf_tsv'
$ws.Range('E50').Value = 'print_zone_flows'
$ws.Range('F50').Value = 'przf'
$ws.Range('G50').Value = 'pri_zf'
$ws.Range('H50').Value = 'timprzf'
$ws.Range('J50').Value = 'ntprzf'
$ws.Range('E51').Value = 'print_zone_flows_xyzt'
$ws.Range('F51').Value = 'przf_xyzt'
$ws.Range('G51').Value = 'pri_zf_xyzt'
$ws.Range('H51').Value = 'timprzf_xyzt'
$ws.Range('J51').Value = 'ntprzf_xyzt'
$ws.Range('E53').Value = 'none/prtichead'
$ws.Range('E33').Value = 'none/print_end_of_period'
$ws.Range('E29').Value = 'none/prt_bc'

# --- View state: active cell / scroll position ---
$ws.Activate()
$ws.Range("E29").Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 2
